$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (dbExcel), shifting dbExcel -> D and WebExcel -> E
$ws.Columns.Item(3).Insert()

# New column C header and value (StatQuery)
$ws.Cells.Item(1,3).Value = "StatQuery"
$ws.Cells.Item(2,3).Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN [''Clinical Trial'',''Transcriptomics'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# Match column width of new column C to column B (both hold long wrapped query text)
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Update the active selection/view to B2
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 2
$ws.Range("B2").Select() | Out-Null
